$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed numeric cells (rows 2-6: refreshed financial figures)
$ws.Range("D2").Value = 5155
$ws.Range("E2").Value = 140
$ws.Range("F2").Value = 140
$ws.Range("G2").Value = 102
$ws.Range("H2").Value = 66
$ws.Range("I2").Value = 66
$ws.Range("K2").Value = 55926
$ws.Range("L2").Value = 49187
$ws.Range("M2").Value = 6739
$ws.Range("N2").Value = 6739
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 1467
$ws.Range("Q2").Value = -5780
$ws.Range("R2").Value = -117
$ws.Range("S2").Value = 5275
$ws.Range("T2").Value = 48
$ws.Range("V2").Value = 11805
$ws.Range("W2").Value = 2.71
$ws.Range("X2").Value = 1.28
$ws.Range("Y2").Value = 0.99
$ws.Range("Z2").Value = 0.13
$ws.Range("AA2").Value = 729.86
$ws.Range("AB2").Value = 359.43
$ws.Range("AC2").Value = 225
$ws.Range("AD2").Value = 45.05
$ws.Range("AE2").Value = 22972
$ws.Range("AF2").Value = 0.44
$ws.Range("AG2").Value = 150
$ws.Range("AH2").Value = 1.48
$ws.Range("AI2").Value = 66.58
$ws.Range("AJ2").Value = 29337111
$ws.Range("D3").Value = 5556
$ws.Range("E3").Value = 682
$ws.Range("F3").Value = 682
$ws.Range("G3").Value = 688
$ws.Range("H3").Value = 504
$ws.Range("I3").Value = 504
$ws.Range("K3").Value = 58526
$ws.Range("L3").Value = 51326
$ws.Range("M3").Value = 7200
$ws.Range("N3").Value = 7200
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 1467
$ws.Range("Q3").Value = -227
$ws.Range("R3").Value = -1425
$ws.Range("S3").Value = 2338
$ws.Range("T3").Value = 64
$ws.Range("V3").Value = 14663
$ws.Range("W3").Value = 12.28
$ws.Range("X3").Value = 9.06
$ws.Range("Y3").Value = 7.22
$ws.Range("Z3").Value = 0.88
$ws.Range("AA3").Value = 712.86
$ws.Range("AB3").Value = 390.85
$ws.Range("AC3").Value = 1716
$ws.Range("AD3").Value = 5.83
$ws.Range("AE3").Value = 24542
$ws.Range("AF3").Value = 0.41
$ws.Range("AG3").Value = 450
$ws.Range("AH3").Value = 4.5
$ws.Range("AI3").Value = 26.22
$ws.Range("AJ3").Value = 29337111
$ws.Range("D4").Value = 5634
$ws.Range("E4").Value = 528
$ws.Range("F4").Value = 528
$ws.Range("G4").Value = 528
$ws.Range("H4").Value = 398
$ws.Range("I4").Value = 398
$ws.Range("K4").Value = 58871
$ws.Range("L4").Value = 50709
$ws.Range("M4").Value = 8162
$ws.Range("N4").Value = 8162
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 1467
$ws.Range("Q4").Value = 1979
$ws.Range("R4").Value = -553
$ws.Range("S4").Value = -1953
$ws.Range("T4").Value = 68
$ws.Range("V4").Value = 12014
$ws.Range("W4").Value = 9.369999999999999
$ws.Range("X4").Value = 7.06
$ws.Range("Y4").Value = 5.18
$ws.Range("Z4").Value = 0.68
$ws.Range("AA4").Value = 621.27
$ws.Range("AB4").Value = 456.43
$ws.Range("AC4").Value = 1356
$ws.Range("AD4").Value = 7
$ws.Range("AE4").Value = 27822
$ws.Range("AF4").Value = 0.34
$ws.Range("AG4").Value = 400
$ws.Range("AH4").Value = 4.21
$ws.Range("AI4").Value = 29.49
$ws.Range("AJ4").Value = 29337111
$ws.Range("D5").Value = 5213
$ws.Range("E5").Value = 668
$ws.Range("F5").Value = 668
$ws.Range("G5").Value = 660
$ws.Range("H5").Value = 502
$ws.Range("I5").Value = 502
$ws.Range("K5").Value = 70252
$ws.Range("L5").Value = 61798
$ws.Range("M5").Value = 8454
$ws.Range("N5").Value = 8454
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 1467
$ws.Range("Q5").Value = -69
$ws.Range("R5").Value = -3848
$ws.Range("S5").Value = 3903
$ws.Range("T5").Value = 67
$ws.Range("V5").Value = 12754
$ws.Range("W5").Value = 12.81
$ws.Range("X5").Value = 9.630000000000001
$ws.Range("Y5").Value = 6.04
$ws.Range("Z5").Value = 0.78
$ws.Range("AA5").Value = 731.03
$ws.Range("AB5").Value = 476.3
$ws.Range("AC5").Value = 1711
$ws.Range("AD5").Value = 6.37
$ws.Range("AE5").Value = 28815
$ws.Range("AF5").Value = 0.38
$ws.Range("AG5").Value = 400
$ws.Range("AH5").Value = 3.67
$ws.Range("AI5").Value = 23.37
$ws.Range("AJ5").Value = 29337111
$ws.Range("D6").Value = 6190
$ws.Range("E6").Value = 681
$ws.Range("F6").Value = 681
$ws.Range("G6").Value = 692
$ws.Range("H6").Value = 506
$ws.Range("I6").Value = 506
$ws.Range("K6").Value = 66864
$ws.Range("L6").Value = 57995
$ws.Range("M6").Value = 8869
$ws.Range("N6").Value = 8869
$ws.Range("P6").Value = 1467
$ws.Range("Q6").Value = -486
$ws.Range("R6").Value = -793
$ws.Range("S6").Value = 2036
$ws.Range("T6").Value = 73
$ws.Range("V6").Value = 15637
$ws.Range("W6").Value = 11.01
$ws.Range("X6").Value = 8.17
$ws.Range("Y6").Value = 5.84
$ws.Range("Z6").Value = 0.74
$ws.Range("AA6").Value = 653.89
$ws.Range("AB6").Value = 504.64
$ws.Range("AC6").Value = 1724
$ws.Range("AD6").Value = 5.01
$ws.Range("AE6").Value = 30232
$ws.Range("AF6").Value = 0.29
$ws.Range("AG6").Value = 450
$ws.Range("AH6").Value = 5.21
$ws.Range("AI6").Value = 26.1
$ws.Range("AJ6").Value = 29337111

# Clear cells removed from the data set (duplicate/erroneous columns, and
# rows 7-9 estimate data that was pulled entirely)
$ws.Range("J2").ClearContents()
$ws.Range("U2").ClearContents()
$ws.Range("J3").ClearContents()
$ws.Range("U3").ClearContents()
$ws.Range("J4").ClearContents()
$ws.Range("U4").ClearContents()
$ws.Range("J5").ClearContents()
$ws.Range("U5").ClearContents()
$ws.Range("U6").ClearContents()
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("P7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("W7").ClearContents()
$ws.Range("X7").ClearContents()
$ws.Range("Y7").ClearContents()
$ws.Range("Z7").ClearContents()
$ws.Range("AA7").ClearContents()
$ws.Range("AC7").ClearContents()
$ws.Range("AD7").ClearContents()
$ws.Range("AE7").ClearContents()
$ws.Range("AF7").ClearContents()
$ws.Range("AG7").ClearContents()
$ws.Range("AH7").ClearContents()
$ws.Range("AI7").ClearContents()
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("P8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("W8").ClearContents()
$ws.Range("X8").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("Z8").ClearContents()
$ws.Range("AA8").ClearContents()
$ws.Range("AC8").ClearContents()
$ws.Range("AD8").ClearContents()
$ws.Range("AE8").ClearContents()
$ws.Range("AF8").ClearContents()
$ws.Range("AG8").ClearContents()
$ws.Range("AH8").ClearContents()
$ws.Range("AI8").ClearContents()
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("P9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("Z9").ClearContents()
$ws.Range("AA9").ClearContents()
$ws.Range("AC9").ClearContents()
$ws.Range("AD9").ClearContents()
$ws.Range("AE9").ClearContents()
$ws.Range("AF9").ClearContents()
$ws.Range("AG9").ClearContents()
$ws.Range("AH9").ClearContents()
$ws.Range("AI9").ClearContents()
